$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.312.84"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "1.865.07"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.69"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4676"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2839"
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06520"
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.73"
$ws.Range("E10").Value = "  +7.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07928"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.43"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "1.873.11"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.152"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6781"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "279.57"
$ws.Range("E16").Value = "  -1.88%  "
$ws.Range("D17").Value = "30.309.27"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.37"
$ws.Range("E18").Value = "  +6.12%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.398"
$ws.Range("E20").Value = "  -1.94%  "
$ws.Range("D21").Value = "2.116.16"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007305"
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.157"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.78"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.167"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.08"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.933"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.389"
$ws.Range("E29").Value = "  +3.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09722"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.396"
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.086"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04739"
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("E35").Value = "  +3.87%  "
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01866"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.577"
$ws.Range("E39").Value = "  +2.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.320"
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.76"
$ws.Range("E41").Value = "  +3.32%  "
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8503"
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4185"
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.30"
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "968.19"
$ws.Range("E47").Value = "  -5.17%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.396"
$ws.Range("E48").Value = "  +3.29%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.190"
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.12"
$ws.Range("E50").Value = "  +0.79%  "
$ws.Range("E51").Value = "  -1.08%  "
